$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 6410525
$ws.Range("I8").Value = 7692430
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 23077290
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = -23077151
$ws.Range("N8").Value = -3278

# Row 28
$ws.Range("H28").Value = 382.16
$ws.Range("I28").Value = 220.47058
$ws.Range("K28").Value = 220.47058
$ws.Range("M28").Value = 264.52942

# Row 105
$ws.Range("H105").Value = 35000
$ws.Range("J105").Value = 35000
$ws.Range("L105").Value = 35000
$ws.Range("N105").Value = -41988

# Row 113
$ws.Range("H113").Value = 3424.3
$ws.Range("J113").Value = 4136.6
$ws.Range("L113").Value = 4136.6
$ws.Range("N113").Value = -10644.6

# Row 132
$ws.Range("H132").Value = 5068.5483
$ws.Range("I132").Value = 4077.2827
$ws.Range("J132").Value = 7918.4375
$ws.Range("K132").Value = 12231.8481
$ws.Range("L132").Value = 23755.3125
$ws.Range("M132").Value = -9701.848100000001
$ws.Range("N132").Value = -28815.3125

# Row 141
$ws.Range("H141").Value = 10946.429
$ws.Range("I141").Value = 2616.111
$ws.Range("J141").Value = 25941
$ws.Range("K141").Value = 7848.333
$ws.Range("L141").Value = 77823
$ws.Range("M141").Value = -2668.333
$ws.Range("N141").Value = -88183

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 10000000
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

# Row 59
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31608

# Row 132
$ws.Range("H132").Value = 2666518.5
$ws.Range("I132").Value = 12120.692
$ws.Range("J132").Value = 3681435.2
$ws.Range("K132").Value = 36362.076
$ws.Range("L132").Value = 11044305.6
$ws.Range("M132").Value = -33832.076
$ws.Range("N132").Value = -11049365.6

$ws = $wb.Worksheets.Item("BSM")
# Row 43
$ws.Range("H43").Value = 230886.67
$ws.Range("J43").Value = 230886.67
$ws.Range("L43").Value = 230886.67
$ws.Range("N43").Value = -231248.67

# Row 75
$ws.Range("H75").Value = 7089.278
$ws.Range("I75").Value = 3528.9092
$ws.Range("K75").Value = 3528.9092
$ws.Range("M75").Value = -2592.9092

# Row 78
$ws.Range("H78").Value = 7089.278
$ws.Range("I78").Value = 3528.9092
$ws.Range("K78").Value = 10586.7276
$ws.Range("M78").Value = -5906.7276

# Row 107
$ws.Range("H107").Value = 1476.5172
$ws.Range("I107").Value = 1445.3478
$ws.Range("J107").Value = 1596
$ws.Range("K107").Value = 1445.3478
$ws.Range("L107").Value = 1596
$ws.Range("M107").Value = 474.6522
$ws.Range("N107").Value = -5436

# Row 109
$ws.Range("H109").Value = 30195
$ws.Range("J109").Value = 30195
$ws.Range("L109").Value = 30195
$ws.Range("N109").Value = -32969

# Row 122
$ws.Range("H122").Value = 40370
$ws.Range("J122").Value = 40370
$ws.Range("L122").Value = 40370
$ws.Range("N122").Value = -50170

# Row 129
$ws.Range("H129").Value = 46633
$ws.Range("J129").Value = 46633
$ws.Range("L129").Value = 46633
$ws.Range("N129").Value = -56633

$ws = $wb.Worksheets.Item("CRP")
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

# Row 31
$ws.Range("H31").Value = 1248.35
$ws.Range("I31").Value = 1002.94666
$ws.Range("J31").Value = 1984.56
$ws.Range("K31").Value = 1002.94666
$ws.Range("L31").Value = 1984.56
$ws.Range("M31").Value = -707.94666
$ws.Range("N31").Value = -2574.56

# Row 34
$ws.Range("H34").Value = 1248.35
$ws.Range("I34").Value = 1002.94666
$ws.Range("J34").Value = 1984.56
$ws.Range("K34").Value = 1002.94666
$ws.Range("L34").Value = 1984.56
$ws.Range("M34").Value = -800.94666
$ws.Range("N34").Value = -2388.56

# Row 86
$ws.Range("H86").Value = 37040704
$ws.Range("I86").Value = 47622020
$ws.Range("J86").Value = 6084.6665
$ws.Range("K86").Value = 47622020
$ws.Range("L86").Value = 6084.6665
$ws.Range("M86").Value = -47620897
$ws.Range("N86").Value = -8330.666499999999

# Row 89
$ws.Range("H89").Value = 37040704
$ws.Range("I89").Value = 47622020
$ws.Range("J89").Value = 6084.6665
$ws.Range("K89").Value = 238110100
$ws.Range("L89").Value = 30423.3325
$ws.Range("M89").Value = -238104484
$ws.Range("N89").Value = -41655.3325

# Row 134
$ws.Range("H134").Value = 2060.6191
$ws.Range("I134").Value = 1636.75
$ws.Range("J134").Value = 2321.4614
$ws.Range("K134").Value = 4910.25
$ws.Range("L134").Value = 6964.3842
$ws.Range("M134").Value = -2375.25
$ws.Range("N134").Value = -12034.3842

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 23809840
$ws.Range("I107").Value = 31250174
$ws.Range("J107").Value = 768.8
$ws.Range("K107").Value = 93750522
$ws.Range("L107").Value = 2306.4
$ws.Range("M107").Value = -93748602
$ws.Range("N107").Value = -6146.4

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1763.1875
$ws.Range("I113").Value = 1611.1
$ws.Range("J113").Value = 2016.6666
$ws.Range("K113").Value = 1611.1
$ws.Range("L113").Value = 2016.6666
$ws.Range("M113").Value = 558.9000000000001
$ws.Range("N113").Value = -6356.6666

# Row 123
$ws.Range("H123").Value = 18918.285
$ws.Range("J123").Value = 18918.285
$ws.Range("L123").Value = 18918.285
$ws.Range("N123").Value = -23818.285

# Row 126
$ws.Range("H126").Value = 4387967.5
$ws.Range("I126").Value = 7577421.5
$ws.Range("J126").Value = 2468.5
$ws.Range("K126").Value = 22732264.5
$ws.Range("L126").Value = 7405.5
$ws.Range("M126").Value = -22729794.5
$ws.Range("N126").Value = -12345.5

# Row 130
$ws.Range("H130").Value = 29620.867
$ws.Range("J130").Value = 37647.7
$ws.Range("L130").Value = 37647.7
$ws.Range("N130").Value = -47687.7

# Row 132
$ws.Range("H132").Value = 3140.5588
$ws.Range("I132").Value = 2027.6923
$ws.Range("J132").Value = 3829.476
$ws.Range("K132").Value = 6083.0769
$ws.Range("L132").Value = 11488.428
$ws.Range("M132").Value = -3553.0769
$ws.Range("N132").Value = -16548.428

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3108.7778
$ws.Range("I61").Value = 3268.8572
$ws.Range("J61").Value = 2548.5
$ws.Range("K61").Value = 3268.8572
$ws.Range("L61").Value = 2548.5
$ws.Range("M61").Value = -3066.8572
$ws.Range("N61").Value = -2952.5

# Row 111
$ws.Range("H111").Value = 34193.5
$ws.Range("J111").Value = 34193.5
$ws.Range("L111").Value = 34193.5
$ws.Range("N111").Value = -42373.5

# Row 113
$ws.Range("H113").Value = 3108.7778
$ws.Range("I113").Value = 3268.8572
$ws.Range("J113").Value = 2548.5
$ws.Range("K113").Value = 3268.8572
$ws.Range("L113").Value = 2548.5
$ws.Range("M113").Value = -1098.8572
$ws.Range("N113").Value = -6888.5

# Row 122
$ws.Range("H122").Value = 8123.476
$ws.Range("I122").Value = 8623.177
$ws.Range("J122").Value = 5999.75
$ws.Range("K122").Value = 25869.531
$ws.Range("L122").Value = 17999.25
$ws.Range("M122").Value = -23419.531
$ws.Range("N122").Value = -22899.25

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1901.5
$ws.Range("I122").Value = 1741.3334
$ws.Range("J122").Value = 2061.6667
$ws.Range("K122").Value = 5224.0002
$ws.Range("L122").Value = 6185.000100000001
$ws.Range("M122").Value = -2774.0002
$ws.Range("N122").Value = -11085.0001
